$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.218.54'

$ws.Range('D3').Value = '1.862.05'
$ws.Range('E3').Value = '  -0.81%  '

$ws.Range('D4').Value = '0.9996'
$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').Value = '0.7141'
$ws.Range('E5').Value = '  -0.55%  '

$ws.Range('D6').Value = '240.58'
$ws.Range('E6').Value = '  +0.14%  '

$ws.Range('E7').Value = '  +0.04%  '

$ws.Range('D8').Value = '0.3087'
$ws.Range('E8').Value = '  -0.44%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07700'
$ws.Range('E9').Value = '  -1.60%  '

$ws.Range('D10').Value = '24.94'
$ws.Range('E10').Value = '  +0.46%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08310'
$ws.Range('E11').Value = '  +0.65%  '

$ws.Range('D12').Value = '1.877.86'
$ws.Range('E12').Value = '  +0.64%  '

$ws.Range('D13').Value = '0.7171'
$ws.Range('E13').Value = '  -1.37%  '

$ws.Range('D14').Value = '5.211'
$ws.Range('E14').Value = '  -1.38%  '

$ws.Range('D15').Value = '90.87'
$ws.Range('E15').Value = '  -0.48%  '

$ws.Range('D16').Value = '29.245.25'
$ws.Range('E16').Value = '  -0.46%  '

$ws.Range('D17').Value = '6.009'
$ws.Range('E17').Value = '  +1.64%  '

$ws.Range('D18').Value = '243.78'
$ws.Range('E18').Value = '  -0.70%  '

$ws.Range('D19').Value = '2.154.50'
$ws.Range('E19').Value = '  +1.44%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007811'
$ws.Range('E20').Value = '  -1.16%  '

$ws.Range('E21').Value = '  -1.08%  '

$ws.Range('D22').Value = '0.9999'
$ws.Range('E22').Value = '  +0.09%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.950'
$ws.Range('E23').Value = '  +0.48%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.000'
$ws.Range('E24').Value = '  +0.04%  '

$ws.Range('D25').Value = '0.1612'
$ws.Range('E25').Value = '  +3.10%  '

$ws.Range('D26').Value = '162.82'
$ws.Range('E26').Value = '  -0.70%  '

$ws.Range('D27').Value = '8.908'
$ws.Range('E27').Value = '  -1.29%  '

$ws.Range('E28').Value = '  +1.37%  '

$ws.Range('D29').Value = '1.354'
$ws.Range('E29').Value = '  -0.57%  '

$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = '1.498'
$ws.Range('E30').Value = '  +0.76%  '

$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '4.443'
$ws.Range('E31').Value = '  +1.06%  '

$ws.Range('D32').Value = '4.259'
$ws.Range('E32').Value = '  +2.82%  '

$ws.Range('D33').Value = '0.05184'
$ws.Range('E33').Value = '  -1.91%  '

$ws.Range('D34').Value = '0.8081'
$ws.Range('E34').Value = '  +11.80%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.930'
$ws.Range('E35').Value = '  -0.28%  '

$ws.Range('D36').Value = '1.173'
$ws.Range('E36').Value = '  -2.41%  '

$ws.Range('D37').Value = '2.686'
$ws.Range('E37').Value = '  +0.36%  '

$ws.Range('E38').Value = '  -0.58%  '

$ws.Range('E39').Value = '  -1.02%  '

$ws.Range('D40').Value = '1.169.42'
$ws.Range('E40').Value = '  -5.29%  '

$ws.Range('D41').Value = '6.231'
$ws.Range('E41').Value = '  +2.35%  '

$ws.Range('D42').Value = '0.9048'
$ws.Range('E42').Value = '  -0.54%  '

$ws.Range('D43').Value = '72.76'
$ws.Range('E43').Value = '  -0.91%  '

$ws.Range('D44').Value = '0.9999'
$ws.Range('E44').Value = '  -0.01%  '

$ws.Range('D45').Value = '2.052.11'
$ws.Range('E45').Value = '  +1.73%  '

$ws.Range('D46').Value = '101.94'
$ws.Range('E46').Value = '  -1.77%  '

$ws.Range('D47').Value = '0.5167'
$ws.Range('E47').Value = '  -3.10%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.780'
$ws.Range('E48').Value = '  +1.29%  '

$ws.Range('D49').Value = '9.373'

$ws.Range('E50').Value = '  -0.94%  '

$ws.Range('E51').Value = '  -0.01%  '
